# Turn C3 from a text "35" into a real numeric 35.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C3").Value = 35

# Add the new row of data (row 4): sudhanshu / Kumar / 25.
# Age "25" must stay a text value (matches the inlineStr in the target),
# so we lead with an apostrophe to force text entry, then reset the
# style back to Normal so we don't leave a stray quote-prefixed format
# behind (the diff shows no style change).
$ws.Range("A4").Value = "sudhanshu "
$ws.Range("B4").Value = "Kumar "
$ws.Range("C4").Value = "'25"
$ws.Range("C4").Style = "Normal"
